$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new data rows above the current row 20 (shifts old rows 20-38
# down to 22-40), making room for two new weekly price observations.
$ws.Rows("20:21").Insert()

# New row 20: Castle Brite / Primera, week of 2021-12-10, Provincia de Los Andes
$ws.Range("A20").Value = 9
$ws.Range("B20").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C20").Value = "Metropolitana"
$ws.Range("D20").Value = 44540
$ws.Range("E20").Value = 13
$ws.Range("F20").Value = "Fruta"
$ws.Range("G20").Value = 100103
$ws.Range("H20").Value = "Frutos de hueso (carozo)"
$ws.Range("I20").Value = 100103003
$ws.Range("J20").Value = "Damasco"
$ws.Range("K20").Value = "Castle Brite"
$ws.Range("L20").Value = "Primera"
$ws.Range("M20").Value = 330
$ws.Range("N20").Value = 16000
$ws.Range("O20").Value = 16000
$ws.Range("P20").Value = 16000
$ws.Range("Q20").Value = "`$/caja 18 kilos"
$ws.Range("R20").Value = "Provincia de Los Andes"
$ws.Range("S20").Value = 889
$ws.Range("T20").Value = 18

# New row 21: Castle Brite / Segunda, week of 2021-12-10, Provincia de Los Andes
$ws.Range("A21").Value = 9
$ws.Range("B21").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C21").Value = "Metropolitana"
$ws.Range("D21").Value = 44540
$ws.Range("E21").Value = 13
$ws.Range("F21").Value = "Fruta"
$ws.Range("G21").Value = 100103
$ws.Range("H21").Value = "Frutos de hueso (carozo)"
$ws.Range("I21").Value = 100103003
$ws.Range("J21").Value = "Damasco"
$ws.Range("K21").Value = "Castle Brite"
$ws.Range("L21").Value = "Segunda"
$ws.Range("M21").Value = 280
$ws.Range("N21").Value = 14000
$ws.Range("O21").Value = 14000
$ws.Range("P21").Value = 14000
$ws.Range("Q21").Value = "`$/caja 18 kilos"
$ws.Range("R21").Value = "Provincia de Los Andes"
$ws.Range("S21").Value = 778
$ws.Range("T21").Value = 18
